# BOT; UPDATE DATA
# Adds one new day (2020-05-07 / serial 43959) of COVID-19 case data to the
# "all", "kobe" and "other" sheets, and bumps the running-total cell on
# sheet "all" (B30) to reflect the corrected cumulative count.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "all": update B30, insert new row 31 (date 43959), footer shifts
# from row 31 to row 32.
# ------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

$wsAll.Range("B30").Value = 273

$wsAll.Rows(31).Insert()
$wsAll.Range("A31").Value = 43959
$wsAll.Range("B31").Value = 273
$wsAll.Range("C31").Value = 268
$wsAll.Range("D31").Value = 79
$wsAll.Range("E31").Value = 69
$wsAll.Range("F31").Value = 10
$wsAll.Range("G31").Value = 8
$wsAll.Range("H31").Value = 181

$wsAll.Range("A31").Select()

# ------------------------------------------------------------------
# Sheet "kobe": update D85/E85, insert new row 86 (date 43959), footer
# shifts from row 86 to row 87.
# ------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("D85").Value = 1
$wsKobe.Range("E85").Value = 273

$wsKobe.Rows(86).Insert()
$wsKobe.Range("A86").Value = 43959
$wsKobe.Range("B86").Value = 0
$wsKobe.Range("C86").Value = 2417
$wsKobe.Range("D86").Value = 0
$wsKobe.Range("E86").Value = 273
$wsKobe.Range("F86").Value = 74
$wsKobe.Range("G86").Value = 65
$wsKobe.Range("H86").Value = 9
$wsKobe.Range("I86").Value = 8
$wsKobe.Range("J86").Value = 174

$wsKobe.Range("A86").Select()

# ------------------------------------------------------------------
# Sheet "other": insert new row 61 (date 43959), footer shifts from row
# 61 to row 62.
# ------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Rows(61).Insert()
$wsOther.Range("A61").Value = 43959
$wsOther.Range("B61").Value = 0
$wsOther.Range("C61").Value = 12
$wsOther.Range("D61").Value = 5
$wsOther.Range("E61").Value = 4
$wsOther.Range("F61").Value = 1
$wsOther.Range("G61").Value = 0
$wsOther.Range("H61").Value = 7

$wsOther.Range("A61").Select()

$wsAll.Activate()
